# Fixed Test data to import.
# Trims the product list down to 4 real rows (was 6), renumbers the ID
# column, swaps several products for their correct category/description/
# image values, and marks the "Images" column as Text-formatted so
# comma-separated filenames don't get mangled on import.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The sheet had 6 data rows (rows 2-7); only 4 are kept. Drop the two
# extra "Thangka" rows (old rows 5 and 6).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Row 2: ID 1 - Buddhas - Aksobhya (category/name/price unchanged, new image)
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "Buddhas - Aksobhya"
$ws.Cells.Item(2,3).Value = "Aksobhya Kupferstatue 21cm"
$ws.Cells.Item(2,4).Value = 649
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = "BAkso008.JPG"
$ws.Cells.Item(2,7).Value = "Beschreibung"

# Row 5 becomes "Thangkas - Big Tsagil" (old row 4 content, ID 4)
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "Thangkas - Big Tsagil"
$ws.Cells.Item(5,3).Value = "große Thangka"
$ws.Cells.Item(5,4).Value = 250
$ws.Cells.Item(5,5).Value = 2

# Row 4 becomes "Malas - Arm Mala" (old row 7 content, ID 3)
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "Malas - Arm Mala"
$ws.Cells.Item(4,3).Value = "kleine mala"
$ws.Cells.Item(4,4).Value = 250
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = "MAM020.JPG"
$ws.Cells.Item(4,7).Value = "Beschreibung"

# Row 3: ID 2 - re-categorised from "Buddhas - Manjusri" to "Buddhas - Shakyamuni"
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "Buddhas - Shakyamuni"
$ws.Cells.Item(3,3).Value = "manjusri"
$ws.Cells.Item(3,4).Value = 649
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = "BSha001.JPG, BSha002.JPG, BSha003.JPG"
$ws.Cells.Item(3,7).Value = "Beschreibung"

$ws.Cells.Item(5,6).Value = "TBT020.JPG, TBT019.JPG, TBT018.JPG"
$ws.Cells.Item(5,7).Value = "Beschreibung"

# Images column must be Text-formatted so multi-file lists like
# "BSha001.JPG, BSha002.JPG, BSha003.JPG" import as literal strings.
$ws.Range("F1:F5").NumberFormat = "@"

$ws.Range("H19").Select()
